$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
# Row 12
$ws.Range("H12").Value = 1599.2307
$ws.Range("I12").Value = 1367
$ws.Range("K12").Value = 1367
$ws.Range("M12").Value = -1197

# Row 33
$ws.Range("H33").Value = 312.1875
$ws.Range("I33").Value = 321.5
$ws.Range("J33").Value = 247
$ws.Range("K33").Value = 321.5
$ws.Range("L33").Value = 247
$ws.Range("M33").Value = -92.5
$ws.Range("N33").Value = -705

# Row 53
$ws.Range("H53").Value = 1507.9
$ws.Range("I53").Value = 258.75
$ws.Range("K53").Value = 258.75
$ws.Range("M53").Value = 378.25

# Row 104
$ws.Range("H104").Value = 626.375
$ws.Range("I104").Value = 513
$ws.Range("K104").Value = 1539
$ws.Range("M104").Value = 208

# Row 111
$ws.Range("H111").Value = 23850.857
$ws.Range("I111").Value = 23850.857
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 71552.571
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -68485.571
$ws.Range("N111").ClearContents()

# Row 137
$ws.Range("H137").Value = 1958.5238
$ws.Range("I137").Value = 2087.3684
$ws.Range("J137").Value = 1852.0869
$ws.Range("K137").Value = 6262.1052
$ws.Range("L137").Value = 5556.2607
$ws.Range("M137").Value = -3712.1052
$ws.Range("N137").Value = -10656.2607

# Row 138
$ws.Range("H138").Value = 223864.44
$ws.Range("I138").Value = 1010.26086
$ws.Range("J138").Value = 456848.38
$ws.Range("K138").Value = 3030.78258
$ws.Range("L138").Value = 1370545.14
$ws.Range("M138").Value = 2109.21742
$ws.Range("N138").Value = -1380825.14

# Row 141
$ws.Range("H141").Value = 2559
$ws.Range("I141").Value = 2015.1818
$ws.Range("K141").Value = 6045.5454
$ws.Range("M141").Value = -865.5454

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 2530.3389
$ws.Range("I32").Value = 2573.5615
$ws.Range("J32").Value = 1298.5
$ws.Range("K32").Value = 2573.5615
$ws.Range("L32").Value = 1298.5
$ws.Range("M32").Value = -2286.5615
$ws.Range("N32").Value = -1872.5

# Row 110
$ws.Range("H110").Value = 4406.6665
$ws.Range("I110").Value = 1644
$ws.Range("K110").Value = 1644
$ws.Range("M110").Value = 401

# Row 132
$ws.Range("H132").Value = 5715.4443
$ws.Range("I132").Value = 2540.3572
$ws.Range("K132").Value = 7621.071599999999
$ws.Range("M132").Value = -5091.071599999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
# Row 20
$ws.Range("H20").Value = 2506.3157
$ws.Range("I20").Value = 2635.4167
$ws.Range("K20").Value = 2635.4167
$ws.Range("M20").Value = -2388.4167

# Row 105
$ws.Range("H105").Value = 2616.125
$ws.Range("I105").Value = 1701.6666
$ws.Range("K105").Value = 1701.6666
$ws.Range("M105").Value = 45.33339999999998

# Row 134
$ws.Range("H134").Value = 2372.9211
$ws.Range("I134").Value = 1691.2051
$ws.Range("K134").Value = 5073.615299999999
$ws.Range("M134").Value = -2538.615299999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
# Row 12
$ws.Range("H12").Value = 957
$ws.Range("J12").Value = 1366.3334
$ws.Range("L12").Value = 1366.3334
$ws.Range("N12").Value = -1706.3334

# Row 31
$ws.Range("H31").Value = 3509.4614
$ws.Range("I31").Value = 1555.8572
$ws.Range("K31").Value = 1555.8572
$ws.Range("M31").Value = -1260.8572

# Row 34
$ws.Range("H34").Value = 3509.4614
$ws.Range("I34").Value = 1555.8572
$ws.Range("K34").Value = 1555.8572
$ws.Range("M34").Value = -1353.8572

# Row 94
$ws.Range("H94").Value = 1636.7
$ws.Range("I94").Value = 983.4
$ws.Range("K94").Value = 983.4
$ws.Range("M94").Value = -532.4

# Row 99
$ws.Range("H99").Value = 4137.276
$ws.Range("I99").Value = 4453.522
$ws.Range("J99").Value = 2925
$ws.Range("K99").Value = 4453.522
$ws.Range("L99").Value = 2925
$ws.Range("M99").Value = -2955.522
$ws.Range("N99").Value = -5921

# Row 126
$ws.Range("H126").Value = 4137.276
$ws.Range("I126").Value = 4453.522
$ws.Range("J126").Value = 2925
$ws.Range("K126").Value = 13360.566
$ws.Range("L126").Value = 8775
$ws.Range("M126").Value = -10890.566
$ws.Range("N126").Value = -13715

# Row 132
$ws.Range("H132").Value = 1775.909

# Row 134
$ws.Range("H134").Value = 1970.3793
$ws.Range("I134").Value = 1909.6538
$ws.Range("J134").Value = 2496.6667
$ws.Range("K134").Value = 5728.9614
$ws.Range("L134").Value = 7490.000100000001
$ws.Range("M134").Value = -3193.9614
$ws.Range("N134").Value = -12560.0001

# Row 135
$ws.Range("H135").Value = 87000
$ws.Range("J135").Value = 87000
$ws.Range("L135").Value = 87000
$ws.Range("N135").Value = -97140

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
# Row 4
$ws.Range("H4").Value = 1860817
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# Row 131
$ws.Range("H131").Value = 1109.5555
$ws.Range("J131").Value = 1693
$ws.Range("L131").Value = 5079
$ws.Range("N131").Value = -15159

# Row 133
$ws.Range("H133").Value = 8860.429
$ws.Range("I133").Value = 8030
$ws.Range("J133").Value = 8998.833000000001
$ws.Range("K133").Value = 24090
$ws.Range("L133").Value = 26996.499
$ws.Range("M133").Value = -19030
$ws.Range("N133").Value = -37116.499

# Row 134
$ws.Range("H134").Value = 3602.25
$ws.Range("I134").Value = 2116.8572
$ws.Range("K134").Value = 6350.571599999999
$ws.Range("M134").Value = -1280.571599999999

# Row 139
$ws.Range("H139").Value = 3391.0334
$ws.Range("I139").Value = 2041
$ws.Range("K139").Value = 6123
$ws.Range("M139").Value = -983

# Row 140
$ws.Range("H140").Value = 45418.5
$ws.Range("I140").Value = 65344.438
$ws.Range("K140").Value = 196033.314
$ws.Range("M140").Value = -190853.314

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
# Row 2
$ws.Range("H2").Value = 89.411766
$ws.Range("I2").Value = 96.61539
$ws.Range("J2").Value = 66
$ws.Range("K2").Value = 96.61539
$ws.Range("L2").Value = 66
$ws.Range("M2").Value = 16.38461
$ws.Range("N2").Value = -292

# Row 3
$ws.Range("H3").Value = 6070157.5
$ws.Range("I3").Value = 2289
$ws.Range("K3").Value = 2289
$ws.Range("M3").Value = -2173

# Row 14
$ws.Range("H14").Value = 5002251
$ws.Range("I14").Value = 5002251
$ws.Range("K14").Value = 5002251
$ws.Range("M14").Value = -5002083

# Row 20
$ws.Range("H20").Value = 22702.6
$ws.Range("I20").Value = 9502.5
$ws.Range("K20").Value = 9502.5
$ws.Range("M20").Value = -9257.5

# Row 21
$ws.Range("H21").Value = 17726.25
$ws.Range("J21").Value = 25454.5
$ws.Range("L21").Value = 25454.5
$ws.Range("N21").Value = -25800.5

# Row 30
$ws.Range("H30").Value = 17726.25
$ws.Range("J30").Value = 25454.5
$ws.Range("L30").Value = 25454.5
$ws.Range("N30").Value = -25664.5

# Row 80
$ws.Range("H80").Value = 8666
$ws.Range("I80").Value = 8666
$ws.Range("K80").Value = 8666
$ws.Range("M80").Value = -7668

# Row 83
$ws.Range("H83").Value = 8666
$ws.Range("I83").Value = 8666
$ws.Range("K83").Value = 43330
$ws.Range("M83").Value = -38338

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
# Row 16
$ws.Range("H16").Value = 920
$ws.Range("I16").Value = 704
$ws.Range("J16").Value = 1352
$ws.Range("K16").Value = 704
$ws.Range("L16").Value = 1352
$ws.Range("M16").Value = -534
$ws.Range("N16").Value = -1692

# Row 22
$ws.Range("H22").Value = 3032.1333
$ws.Range("I22").Value = 783.2857
$ws.Range("J22").Value = 4999.875
$ws.Range("K22").Value = 783.2857
$ws.Range("L22").Value = 4999.875
$ws.Range("M22").Value = -488.2857
$ws.Range("N22").Value = -5589.875

# Row 27
$ws.Range("H27").Value = 3032.1333
$ws.Range("I27").Value = 783.2857
$ws.Range("J27").Value = 4999.875
$ws.Range("K27").Value = 783.2857
$ws.Range("L27").Value = 4999.875
$ws.Range("M27").Value = -676.2857
$ws.Range("N27").Value = -5213.875

# Row 31
$ws.Range("H31").Value = 4167.25
$ws.Range("I31").Value = 1400
$ws.Range("K31").Value = 1400
$ws.Range("M31").Value = -1152

# Row 40
$ws.Range("H40").Value = 5673.5415
$ws.Range("I40").Value = 5050.6523
$ws.Range("K40").Value = 5050.6523
$ws.Range("M40").Value = -4914.6523

# Row 63
$ws.Range("H63").Value = 91883
$ws.Range("J63").Value = 91883
$ws.Range("L63").Value = 91883
$ws.Range("N63").Value = -93381

# Row 66
$ws.Range("H66").Value = 91883
$ws.Range("J66").Value = 91883
$ws.Range("L66").Value = 275649
$ws.Range("N66").Value = -283137

# Row 109
$ws.Range("H109").Value = 98995
$ws.Range("J109").Value = 98995
$ws.Range("L109").Value = 98995
$ws.Range("N109").Value = -101769

# Row 132
$ws.Range("H132").Value = 3767.6667
$ws.Range("I132").Value = 1526.6364
$ws.Range("K132").Value = 4579.9092
$ws.Range("M132").Value = -2049.9092

# Row 136
$ws.Range("H136").Value = 4644.981
$ws.Range("I136").Value = 4064.681
$ws.Range("K136").Value = 12194.043
$ws.Range("M136").Value = -9644.043

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
# Row 21
$ws.Range("H21").Value = 15996.667
$ws.Range("I21").Value = 14507.5
$ws.Range("J21").Value = 18975
$ws.Range("K21").Value = 14507.5
$ws.Range("L21").Value = 18975
$ws.Range("M21").Value = -14272.5
$ws.Range("N21").Value = -19445

# Row 32
$ws.Range("H32").Value = 6900
$ws.Range("I32").Value = 6900
$ws.Range("K32").Value = 6900
$ws.Range("M32").Value = -6583

# Row 35
$ws.Range("H35").Value = 15996.667
$ws.Range("I35").Value = 14507.5
$ws.Range("J35").Value = 18975
$ws.Range("K35").Value = 14507.5
$ws.Range("L35").Value = 18975
$ws.Range("M35").Value = -14217.5
$ws.Range("N35").Value = -19555

# Row 139
$ws.Range("H139").Value = 78943
$ws.Range("I139").Value = 60000
$ws.Range("J139").Value = 83678.75
$ws.Range("K139").Value = 60000
$ws.Range("L139").Value = 83678.75
$ws.Range("M139").Value = -54860
$ws.Range("N139").Value = -93958.75
